# Generate Report for Handback
# Refresh the handoff/handback timestamps recorded for the
# "048ba33b-3a53-4b93-8c8f-5980f837820b.md" file after a new handback
# report run, on both locale sheets, and roll the newest timestamp up
# into the Overview sheet's "Latest HO Xliff Generate Date" column.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: row 2 is the 048ba33b-... file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-10-14 07:54:53"
$wsZhCn.Range("K2").Value = "2016-10-14 07:55:35"

# --- de-de sheet: row 2 is the 048ba33b-... file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-10-14 07:55:04"
$wsDeDe.Range("K2").Value = "2016-10-14 07:55:52"

# --- Overview sheet: roll up the latest HO xliff generate date ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-10-14 07:55:04"
